# Auto-generated edit script: applies the 2024-05-10 daily violent-crime
# data refresh to the "2024" (column K) figures across all affected sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 2596
$ws.Range("K3").Value = 2507
$ws.Range("K4").Value = 524
$ws.Range("K5").Value = 164
$ws.Range("K6").Value = 3128
$ws.Range("K7").Value = 8919

$ws = $wb.Worksheets.Item('Norwood Park')
$ws.Range("K6").Value = 4
$ws.Range("K7").Value = 23

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("K3").Value = 29
$ws.Range("K7").Value = 131

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K2").Value = 173
$ws.Range("K3").Value = 177
$ws.Range("K6").Value = 198
$ws.Range("K7").Value = 594

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("K6").Value = 98
$ws.Range("K7").Value = 351

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("K3").Value = 47
$ws.Range("K6").Value = 40
$ws.Range("K7").Value = 142

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K3").Value = 100
$ws.Range("K6").Value = 90
$ws.Range("K7").Value = 290

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("K3").Value = 60
$ws.Range("K7").Value = 162

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K7").Value = 266
$ws.Range("K8").Value = 594
$ws.Range("K10").Value = 50
$ws.Range("K11").Value = 189
$ws.Range("K14").Value = 53
$ws.Range("K17").Value = 15
$ws.Range("K18").Value = 60
$ws.Range("K19").Value = 262
$ws.Range("K20").Value = 203
$ws.Range("K27").Value = 96
$ws.Range("K29").Value = 457
$ws.Range("K31").Value = 104
$ws.Range("K33").Value = 351
$ws.Range("K35").Value = 14
$ws.Range("K36").Value = 104
$ws.Range("K37").Value = 290
$ws.Range("K42").Value = 308
$ws.Range("K48").Value = 111
$ws.Range("K50").Value = 58
$ws.Range("K52").Value = 244
$ws.Range("K53").Value = 131
$ws.Range("K54").Value = 165
$ws.Range("K55").Value = 96
$ws.Range("K63").Value = 37
$ws.Range("K66").Value = 31
$ws.Range("K67").Value = 349
$ws.Range("K69").Value = 23
$ws.Range("K71").Value = 25
$ws.Range("K76").Value = 130
$ws.Range("K78").Value = 124
$ws.Range("K85").Value = 428
$ws.Range("K86").Value = 56
$ws.Range("K88").Value = 103
$ws.Range("K89").Value = 117
$ws.Range("K90").Value = 79
$ws.Range("K91").Value = 84
$ws.Range("K95").Value = 142
$ws.Range("K96").Value = 126
$ws.Range("K97").Value = 76
$ws.Range("K99").Value = 162
$ws.Range("K101").Value = 8919

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("K2").Value = 35
$ws.Range("K7").Value = 104

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("K5").Value = 5
$ws.Range("K7").Value = 349

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("K6").Value = 72
$ws.Range("K7").Value = 165

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K2").Value = 124
$ws.Range("K3").Value = 154
$ws.Range("K7").Value = 457

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("K6").Value = 58
$ws.Range("K7").Value = 111

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("K2").Value = 83
$ws.Range("K6").Value = 90
$ws.Range("K7").Value = 262

$ws = $wb.Worksheets.Item('River North')
$ws.Range("K6").Value = 79
$ws.Range("K7").Value = 130

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("K5").Value = 1
$ws.Range("K7").Value = 53

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("K3").Value = 97
$ws.Range("K7").Value = 308

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("K3").Value = 8
$ws.Range("K7").Value = 50

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("K6").Value = 46
$ws.Range("K7").Value = 124

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("K3").Value = 22
$ws.Range("K7").Value = 96

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("K3").Value = 18
$ws.Range("K6").Value = 59
$ws.Range("K7").Value = 126

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("K6").Value = 19
$ws.Range("K7").Value = 84

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("K3").Value = 57
$ws.Range("K6").Value = 70
$ws.Range("K7").Value = 203

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("K6").Value = 14
$ws.Range("K7").Value = 60

$ws = $wb.Worksheets.Item('Burnside')
$ws.Range("K3").Value = 3
$ws.Range("K6").Value = 4
$ws.Range("K7").Value = 15

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("K6").Value = 23
$ws.Range("K7").Value = 104

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("K2").Value = 86
$ws.Range("K3").Value = 84
$ws.Range("K6").Value = 74
$ws.Range("K7").Value = 266

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("K3").Value = 7
$ws.Range("K6").Value = 35
$ws.Range("K7").Value = 58

$ws = $wb.Worksheets.Item('North Center')
$ws.Range("K2").Value = 9
$ws.Range("K7").Value = 31

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("K6").Value = 75
$ws.Range("K7").Value = 189

$ws = $wb.Worksheets.Item('Gold Coast')
$ws.Range("K6").Value = 11
$ws.Range("K7").Value = 14

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("K6").Value = 48
$ws.Range("K7").Value = 76

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("K2").Value = 22
$ws.Range("K6").Value = 53
$ws.Range("K7").Value = 103

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("K4").Value = 17
$ws.Range("K7").Value = 117

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("K3").Value = 20
$ws.Range("K7").Value = 96

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("K2").Value = 13
$ws.Range("K7").Value = 56

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("K2").Value = 32
$ws.Range("K3").Value = 21
$ws.Range("K7").Value = 79

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("K2").Value = 158
$ws.Range("K3").Value = 147
$ws.Range("K6").Value = 96
$ws.Range("K7").Value = 428

$ws = $wb.Worksheets.Item('Oakland')
$ws.Range("K6").Value = 6
$ws.Range("K7").Value = 25

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("K2").Value = 68
$ws.Range("K7").Value = 244
